$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K9").Value = 169520.3

$ws.Range("M10").Value = 431099.65
$ws.Range("N10").Value = 125812.36
$ws.Range("O10").Value = 77310.06

$ws.Range("K11").Value = 134197.7

$ws.Range("O16").Value = 8408.53

$ws.Range("M17").Value = 918646.48
$ws.Range("N17").Value = 448735.6
$ws.Range("O17").Value = 438712.76

$ws.Range("O18").Value = 488726.58

$ws.Range("N23").Value = 108076
$ws.Range("N24").Value = 5667.92
$ws.Range("N27").Value = 3030.6
